# Update code for report co so
# Adds two new worksheets ("Đơn 1 bác sĩ" and "Lương") after the existing
# "Đơn sale chính" sheet, and populates "Đơn 1 bác sĩ" with a single
# doctor-commission report (header row, one data row, one totals row).

$wb = $excel.ActiveWorkbook

# --- Existing sheet stays first; new sheets are appended after it -------
$mainSheet = $wb.Worksheets.Item(1)

$donSheet = $wb.Worksheets.Add($null, $mainSheet)
$donSheet.Name = "Đơn 1 bác sĩ"

$luongSheet = $wb.Worksheets.Add($null, $donSheet)
$luongSheet.Name = "Lương"

# --- Populate "Đơn 1 bác sĩ" ---------------------------------------------
$ws = $donSheet

# Header row (row 1)
$headers = @(
    "Tiền tố",
    "Mã dịch vụ",
    "Ngày thực hiện",
    "Cơ sở",
    "Khách hàng",
    "Nguồn khách",
    "Tên dịch vụ",
    "Sale chính",
    "Đơn giá gốc",
    "Sale phụ",
    "Upsale",
    "Đơn giá",
    "Thanh toán lần đầu",
    "Trả sau",
    "Đã thanh toán",
    "Dư nợ",
    "Bác sĩ 1",
    "Bác sĩ 2",
    "Phụ phẫu 1",
    "Phụ phẫu 2",
    "Công phụ phẫu 1",
    "Công phụ phẫu 2"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data row (row 2)
$ws.Cells.Item(2, 1).Value = "HD-LUXURY"
$ws.Cells.Item(2, 2).Value = 521
# Force the execution-date column to stay plain text (not an Excel date
# serial) since the source report stores it as "dd-mm-yyyy" text.
$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "07-05-2024"
$ws.Cells.Item(2, 4).Value = "SÓC TRĂNG"
$ws.Cells.Item(2, 5).Value = "đa ni "
$ws.Cells.Item(2, 6).Value = "CTV"
$ws.Cells.Item(2, 7).Value = "Phun môi"
$ws.Cells.Item(2, 8).Value = "Thạch Hoàng Nhân"
$ws.Cells.Item(2, 9).Value = 5500000
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 5500000
$ws.Cells.Item(2, 13).Value = 5500000
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 5500000
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = "Bác Sĩ Ngoài"
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 0
$ws.Cells.Item(2, 22).Value = 0

# Totals row (row 3). Columns with no meaningful total (C, D, E, F, G, H,
# J, Q, R, S, T) are left blank, same as the source report.
$ws.Cells.Item(3, 1).Value = "Tổng"
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 9).Value = 5500000
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 5500000
$ws.Cells.Item(3, 13).Value = 5500000
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 5500000
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 21).Value = 0
$ws.Cells.Item(3, 22).Value = 0

# Leave "Lương" sheet empty (placeholder sheet for future payroll data).

# Select the first sheet to match the original workbook's active state.
$mainSheet.Activate()
